# Apply the "May 9th" changes:
#  - Remove the first three data rows (old rows 2-4), shifting the remaining
#    data rows up.
#  - Append thirteen brand-new data rows after the existing data (new rows
#    19-31), extending the used range to A1:C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first three samples (old A2:C4) - this shifts rows 5-21 up to 2-18.
$ws.Rows("2:4").Delete()

# Append the new samples collected on May 9th.
$newRows = @(
    @(-62.44527816772461, -7.982970237731934, -32.66093826293945),
    @(23.0417537689209, 16.5068416595459, -16.63084411621094),
    @(-5.877625465393066, -9.525984764099119, -0.1709146499633789),
    @(-0.5468623042106628, -5.408160209655762, 43.39143753051758),
    @(-6.48110294342041, 21.88984298706055, 9.76063346862793),
    @(19.5135498046875, -70.08018493652344, 22.09575080871582),
    @(-37.2692756652832, 22.71374130249023, -6.724684715270996),
    @(-22.65700340270996, 4.547637939453125, -17.0935115814209),
    @(78.31330871582031, -70.41346740722656, -6.543253421783447),
    @(-3.919276475906372, 13.1894645690918, 2.611269950866699),
    @(41.53395843505859, -10.51663017272949, 21.93678855895996),
    @(19.39011192321777, -10.6725959777832, -12.10480403900146),
    @(10.54305171966553, 18.1186752319336, -10.15065765380859)
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
}
